$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Sheet1: insert a new leading "id" column with row numbers, keep
#    Name/Age but update Peter's Age value.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$null = $ws1.Columns.Item(1).Insert()
$ws1.Cells.Item(1,1).Value = "id"
$ws1.Cells.Item(2,1).Value = 1
$ws1.Cells.Item(3,1).Value = 2
$ws1.Cells.Item(2,3).Value = 45
$null = $ws1.Range("C2").Select()

# ------------------------------------------------------------------
# 2) Drop Sheet3 entirely - its Id/Date table is no longer needed.
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$null = $ws3.Delete()

# ------------------------------------------------------------------
# 3) Sheet2: replace the old Book/Price sample rows with the new
#    book/price/location table.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Column B holds the price strings - format the column as Text first
# so the "$9.9"/"$12.9" strings are kept verbatim instead of being
# parsed into numbers.
$ws2.Columns.Item(2).NumberFormat = "@"

$ws2.Cells.Item(1,1).Value = "Book"
$ws2.Cells.Item(1,2).Value = "Price"
$ws2.Cells.Item(1,3).Value = "Location"

$ws2.Cells.Item(2,1).Value = "To kill a mocking bird"
$ws2.Cells.Item(2,2).Value = "$9.9"
$ws2.Cells.Item(2,3).Value = "Level-2-rack-1"

$ws2.Cells.Item(3,1).Value = "Python cookbook"
$ws2.Cells.Item(3,2).Value = "$12.9"
$ws2.Cells.Item(3,3).Value = "Level-1-rack-2"

# Column widths to fit the longer book titles / rack locations.
$ws2.Columns.Item(1).ColumnWidth = 20
$ws2.Columns.Item(3).ColumnWidth = 16

# Smaller font for the long title in A2.
$ws2.Cells.Item(2,1).Font.Size = 10.5

# Make Sheet2 the active/selected sheet, matching the new tab order.
$null = $ws2.Select()
$null = $ws2.Range("D8").Select()
